# Feria Lagunitas de Puerto Montt - Perejil: add a new weekly price record.
# A new row is inserted right before the current row 191 (pushing the
# existing rows 191-288 down to 192-289) with a fresh observation dated
# 44813 (sharing the same market/product/region metadata pattern used
# throughout the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 191; this shifts rows 191..288 down to 192..289
# and the new blank row inherits formatting from the row above (so the date
# column keeps its numeric date style).
$ws.Rows(191).Insert()

$ws.Cells.Item(191, 1).Value = 4
$ws.Cells.Item(191, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(191, 3).Value = "Los Lagos"
$ws.Cells.Item(191, 4).Value = 44813
$ws.Cells.Item(191, 5).Value = 10
$ws.Cells.Item(191, 6).Value = 100112044
$ws.Cells.Item(191, 7).Value = "Perejil"
$ws.Cells.Item(191, 8).Value = "Sin especificar"
$ws.Cells.Item(191, 9).Value = "Primera"
$ws.Cells.Item(191, 10).Value = 180
$ws.Cells.Item(191, 11).Value = 6000
$ws.Cells.Item(191, 12).Value = 6000
$ws.Cells.Item(191, 13).Value = 6000
$ws.Cells.Item(191, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(191, 15).Value = "Región Metropolitana"
$ws.Cells.Item(191, 16).Value = 2000
$ws.Cells.Item(191, 17).Value = 3
$ws.Cells.Item(191, 18).Value = "Hortaliza"
